$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at K (everything from K..O shifts right to L..P).
$ws.Range("K1").EntireColumn.Insert()

# The insert leaves a couple of phantom blank/styled cells behind in rows
# that had a cell in the column just to the left (J) but none at K before
# the shift — clear those so the new column stays genuinely empty there,
# matching the source rows that never had a K cell.
$ws.Range("K2").Clear() | Out-Null
$ws.Range("K4").Clear() | Out-Null

# New header cell K1: "Verbose" + newline + "Comments", styled like F1
# (center/center/wrap — same look as the other wrapped header cell).
$ws.Range("K1").Value = "Verbose`nComments"
$ws.Range("K1").HorizontalAlignment = -4108
$ws.Range("K1").VerticalAlignment = -4108
$ws.Range("K1").WrapText = $true

# Move the selection like the author's workbook shows post-edit.
$ws.Range("K2").Select() | Out-Null
